# Auto update Excel log
# 1) ALERTS sheet: append a new FALL_DETECTED / CRITICAL row (row 4)
# 2) mmWave sheet: append nine new PRESENCE_DETECTED / Active rows (rows 27-35)

$wb = $excel.ActiveWorkbook

function Add-LogRow {
    param($ws, $Row, $Date, $Timestamp, $Hour, $Location, $Value, $Status)

    # Force column A to Text format first so the date-shaped string
    # ("2026-01-31") is stored verbatim instead of being auto-converted
    # into a date serial number by Excel's smart input parsing.
    $ws.Range("A$Row").NumberFormat = "@"
    $ws.Range("A$Row").Value = $Date
    $ws.Range("B$Row").Value = $Timestamp
    $ws.Range("C$Row").Value = $Hour
    $ws.Range("D$Row").Value = $Location
    $ws.Range("E$Row").Value = $Value
    $ws.Range("F$Row").Value = $Status
    # Drop back to the default "Normal" style so the new cells don't end
    # up carrying an extra text-format style that the original rows lack.
    $ws.Range("A$Row").Style = "Normal"
}

# --- ALERTS sheet: new row 4 -------------------------------------------------
$alerts = $wb.Worksheets.Item("ALERTS")
Add-LogRow $alerts 4 "2026-01-31" "21:43:40" "21:00" "Living Room" "CRITICAL" "FALL_DETECTED"

# --- mmWave sheet: new rows 27-35 -------------------------------------------
$mmwave = $wb.Worksheets.Item("mmWave")
$timestamps = @(
    "21:43:34",
    "21:43:35",
    "21:43:35",
    "21:43:41",
    "21:43:52",
    "21:44:02",
    "21:44:13",
    "21:44:23",
    "21:44:34"
)

$row = 27
foreach ($ts in $timestamps) {
    Add-LogRow $mmwave $row "2026-01-31" $ts "21:00" "Living Room" "PRESENCE_DETECTED" "Active"
    $row++
}
